$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testcase")

# comparetype: s2tcompare -> likeobjectcompare
$ws.Range("B2").Value = "likeobjectcompare"

# s2tpath: clear value
$ws.Range("B30").Value = ""

# s2tmappingsheet: clear value
$ws.Range("B31").Value = ""

# Update the selection / view to match final state
$ws.Range("C9").Select()
